# Disaggregation of commodity Copper
#
# 1) Rename the shared label "Copper ores and concentrates" -> "Copper".
#    Every yearly sheet (2000..2100) carries this label in cell C4, so we
#    rewrite it on every worksheet.
# 2) A handful of sheets also carry a one-ULP recalculation drift on the
#    cached total in D4 (last-significant-digit change only).

$wb = $excel.ActiveWorkbook

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $ws.Range("C4").Value = "Copper"
}

# Tiny (last-digit) recalculated total updates on specific year sheets.
$wb.Worksheets.Item("2033").Range("D4").Value = 110830.1039065614
$wb.Worksheets.Item("2039").Range("D4").Value = 216811.3829355027
$wb.Worksheets.Item("2041").Range("D4").Value = 278380.1093116245
$wb.Worksheets.Item("2045").Range("D4").Value = 642552.158481146
$wb.Worksheets.Item("2067").Range("D4").Value = 748329.7765664503
$wb.Worksheets.Item("2069").Range("D4").Value = 939284.8480597934
$wb.Worksheets.Item("2072").Range("D4").Value = 1396816.716286596
